$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(15, 8).Value = 4562.293  # H15: 2998.8523 -> 4562.293
$ws.Cells.Item(15, 9).Value = 4562.293  # I15: 2998.8523 -> 4562.293
$ws.Cells.Item(15, 11).Value = 13686.879  # K15: 8996.5569 -> 13686.879
$ws.Cells.Item(15, 13).Value = -13517.879  # M15: -8827.5569 -> -13517.879
$ws.Cells.Item(19, 8).Value = 1368.3182  # H19: 1433.4762 -> 1368.3182
$ws.Cells.Item(19, 9).Value = 1025.25  # I19: 1091.1818 -> 1025.25
$ws.Cells.Item(19, 10).Value = 1780  # J19: 1810 -> 1780
$ws.Cells.Item(19, 11).Value = 1025.25  # K19: 1091.1818 -> 1025.25
$ws.Cells.Item(19, 12).Value = 1780  # L19: 1810 -> 1780
$ws.Cells.Item(19, 13).Value = -850.25  # M19: -916.1818000000001 -> -850.25
$ws.Cells.Item(19, 14).Value = -2130  # N19: -2160 -> -2130
$ws.Cells.Item(28, 8).Value = 1238.5758  # H28: 831.89655 -> 1238.5758
$ws.Cells.Item(28, 9).Value = 644.2593000000001  # I28: 601.7308 -> 644.2593000000001
$ws.Cells.Item(28, 10).Value = 3913  # J28: 2826.6667 -> 3913
$ws.Cells.Item(28, 11).Value = 644.2593000000001  # K28: 601.7308 -> 644.2593000000001
$ws.Cells.Item(28, 12).Value = 3913  # L28: 2826.6667 -> 3913
$ws.Cells.Item(28, 13).Value = -159.2593000000001  # M28: -116.7308 -> -159.2593000000001
$ws.Cells.Item(28, 14).Value = -4883  # N28: -3796.6667 -> -4883
$ws.Cells.Item(107, 8).Value = 451.95834  # H107: 620.95 -> 451.95834
$ws.Cells.Item(107, 9).Value = 451.95834  # I107: 647.6923 -> 451.95834
$ws.Cells.Item(107, 10).Value = 0  # J107: 571.2857 -> 0
$ws.Cells.Item(107, 11).Value = 451.95834  # K107: 647.6923 -> 451.95834
$ws.Cells.Item(107, 12).Value = 0  # L107: 571.2857 -> 0
$ws.Cells.Item(107, 13).Value = 1468.04166  # M107: 1272.3077 -> 1468.04166
$ws.Cells.Item(107, 14).ClearContents()  # N107: -4411.2857 -> (removed)
$ws.Cells.Item(111, 8).Value = 3353.65  # H111: 58766 -> 3353.65
$ws.Cells.Item(111, 9).Value = 2881.6155  # I111: 10000 -> 2881.6155
$ws.Cells.Item(111, 10).Value = 4230.2856  # J111: 107532 -> 4230.2856
$ws.Cells.Item(111, 11).Value = 8644.8465  # K111: 30000 -> 8644.8465
$ws.Cells.Item(111, 12).Value = 12690.8568  # L111: 322596 -> 12690.8568
$ws.Cells.Item(111, 13).Value = -5577.8465  # M111: -26933 -> -5577.8465
$ws.Cells.Item(111, 14).Value = -18824.8568  # N111: -328730 -> -18824.8568
$ws.Cells.Item(116, 8).Value = 2608.077  # H116: 2791.4546 -> 2608.077
$ws.Cells.Item(116, 10).Value = 3015  # J116: 3581.2 -> 3015
$ws.Cells.Item(116, 12).Value = 3015  # L116: 3581.2 -> 3015
$ws.Cells.Item(116, 14).Value = -9899  # N116: -10465.2 -> -9899

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 5982.2466  # H32: 4895.1753 -> 5982.2466
$ws.Cells.Item(32, 9).Value = 3695.7078  # I32: 2960.2354 -> 3695.7078
$ws.Cells.Item(32, 10).Value = 18367.666  # J32: 18601 -> 18367.666
$ws.Cells.Item(32, 11).Value = 3695.7078  # K32: 2960.2354 -> 3695.7078
$ws.Cells.Item(32, 12).Value = 18367.666  # L32: 18601 -> 18367.666
$ws.Cells.Item(32, 13).Value = -3408.7078  # M32: -2673.2354 -> -3408.7078
$ws.Cells.Item(32, 14).Value = -18941.666  # N32: -19175 -> -18941.666
$ws.Cells.Item(63, 8).Value = 3261.111  # H63: 3392.8572 -> 3261.111
$ws.Cells.Item(63, 9).Value = 2840  # I63: 0 -> 2840
$ws.Cells.Item(63, 10).Value = 3423.077  # J63: 3392.8572 -> 3423.077
$ws.Cells.Item(63, 11).Value = 2840  # K63: 0 -> 2840
$ws.Cells.Item(63, 12).Value = 3423.077  # L63: 3392.8572 -> 3423.077
$ws.Cells.Item(63, 13).Value = -2154  # M63: None -> -2154
$ws.Cells.Item(63, 14).Value = -4795.077  # N63: -4764.8572 -> -4795.077
$ws.Cells.Item(66, 8).Value = 3261.111  # H66: 3392.8572 -> 3261.111
$ws.Cells.Item(66, 9).Value = 2840  # I66: 0 -> 2840
$ws.Cells.Item(66, 10).Value = 3423.077  # J66: 3392.8572 -> 3423.077
$ws.Cells.Item(66, 11).Value = 14200  # K66: 0 -> 14200
$ws.Cells.Item(66, 12).Value = 17115.385  # L66: 16964.286 -> 17115.385
$ws.Cells.Item(66, 13).Value = -10768  # M66: None -> -10768
$ws.Cells.Item(66, 14).Value = -23979.385  # N66: -23828.286 -> -23979.385
$ws.Cells.Item(74, 8).Value = 22527.54  # H74: 25396.479 -> 22527.54
$ws.Cells.Item(74, 9).Value = 31394.787  # I74: 39685.92 -> 31394.787
$ws.Cells.Item(74, 10).Value = 7126.5264  # J74: 6820.2 -> 7126.5264
$ws.Cells.Item(74, 11).Value = 31394.787  # K74: 39685.92 -> 31394.787
$ws.Cells.Item(74, 12).Value = 7126.5264  # L74: 6820.2 -> 7126.5264
$ws.Cells.Item(74, 13).Value = -30520.787  # M74: -38811.92 -> -30520.787
$ws.Cells.Item(74, 14).Value = -8874.526399999999  # N74: -8568.200000000001 -> -8874.526399999999
$ws.Cells.Item(77, 8).Value = 22527.54  # H77: 25396.479 -> 22527.54
$ws.Cells.Item(77, 9).Value = 31394.787  # I77: 39685.92 -> 31394.787
$ws.Cells.Item(77, 10).Value = 7126.5264  # J77: 6820.2 -> 7126.5264
$ws.Cells.Item(77, 11).Value = 156973.935  # K77: 198429.6 -> 156973.935
$ws.Cells.Item(77, 12).Value = 35632.632  # L77: 34101 -> 35632.632
$ws.Cells.Item(77, 13).Value = -152605.935  # M77: -194061.6 -> -152605.935
$ws.Cells.Item(77, 14).Value = -44368.632  # N77: -42837 -> -44368.632

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(50, 8).Value = 0  # H50: 30000 -> 0
$ws.Cells.Item(50, 10).Value = 0  # J50: 30000 -> 0
$ws.Cells.Item(50, 12).Value = 0  # L50: 30000 -> 0
$ws.Cells.Item(50, 14).ClearContents()  # N50: -31148 -> (removed)
$ws.Cells.Item(51, 8).Value = 29699.5  # H51: 29733 -> 29699.5
$ws.Cells.Item(51, 10).Value = 29699.5  # J51: 29733 -> 29699.5
$ws.Cells.Item(51, 12).Value = 29699.5  # L51: 29733 -> 29699.5
$ws.Cells.Item(51, 14).Value = -30681.5  # N51: -30715 -> -30681.5
$ws.Cells.Item(55, 8).Value = 29169.5  # H55: 29115.6 -> 29169.5
$ws.Cells.Item(55, 10).Value = 29169.5  # J55: 29115.6 -> 29169.5
$ws.Cells.Item(55, 12).Value = 29169.5  # L55: 29115.6 -> 29169.5
$ws.Cells.Item(55, 14).Value = -29715.5  # N55: -29661.6 -> -29715.5
$ws.Cells.Item(76, 8).Value = 29900  # H76: 30000 -> 29900
$ws.Cells.Item(76, 10).Value = 29900  # J76: 30000 -> 29900
$ws.Cells.Item(76, 12).Value = 29900  # L76: 30000 -> 29900
$ws.Cells.Item(76, 14).Value = -30530  # N76: -30630 -> -30530
$ws.Cells.Item(79, 8).Value = 29900  # H79: 30000 -> 29900
$ws.Cells.Item(79, 10).Value = 29900  # J79: 30000 -> 29900
$ws.Cells.Item(79, 12).Value = 29900  # L79: 30000 -> 29900
$ws.Cells.Item(79, 14).Value = -32084  # N79: -32184 -> -32084
$ws.Cells.Item(99, 8).Value = 1393.6364  # H99: 0 -> 1393.6364
$ws.Cells.Item(99, 9).Value = 1283  # I99: 0 -> 1283
$ws.Cells.Item(99, 10).Value = 2500  # J99: 0 -> 2500
$ws.Cells.Item(99, 11).Value = 1283  # K99: 0 -> 1283
$ws.Cells.Item(99, 12).Value = 2500  # L99: 0 -> 2500
$ws.Cells.Item(99, 13).Value = 215  # M99: None -> 215
$ws.Cells.Item(99, 14).Value = -5496  # N99: None -> -5496
$ws.Cells.Item(134, 8).Value = 731411.4  # H134: 479063.6 -> 731411.4
$ws.Cells.Item(134, 9).Value = 1216480.1  # I134: 692326.9399999999 -> 1216480.1
$ws.Cells.Item(134, 10).Value = 3808.1365  # J134: 3322.2693 -> 3808.1365
$ws.Cells.Item(134, 11).Value = 3649440.3  # K134: 2076980.82 -> 3649440.3
$ws.Cells.Item(134, 12).Value = 11424.4095  # L134: 9966.8079 -> 11424.4095
$ws.Cells.Item(134, 13).Value = -3646905.3  # M134: -2074445.82 -> -3646905.3
$ws.Cells.Item(134, 14).Value = -16494.4095  # N134: -15036.8079 -> -16494.4095

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(87, 8).Value = 0  # H87: 21666.666 -> 0
$ws.Cells.Item(87, 10).Value = 0  # J87: 21666.666 -> 0
$ws.Cells.Item(87, 12).Value = 0  # L87: 21666.666 -> 0
$ws.Cells.Item(87, 14).ClearContents()  # N87: -24038.666 -> (removed)
$ws.Cells.Item(90, 8).Value = 0  # H90: 21666.666 -> 0
$ws.Cells.Item(90, 10).Value = 0  # J90: 21666.666 -> 0
$ws.Cells.Item(90, 12).Value = 0  # L90: 64999.99800000001 -> 0
$ws.Cells.Item(90, 14).ClearContents()  # N90: -76855.99800000001 -> (removed)

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(2, 8).Value = 303421.88  # H2: 303423.8 -> 303421.88
$ws.Cells.Item(2, 9).Value = 855.5833  # I2: 932.63635 -> 855.5833
$ws.Cells.Item(2, 10).Value = 516998.06  # J2: 488279.56 -> 516998.06
$ws.Cells.Item(2, 11).Value = 5133.4998  # K2: 5595.8181 -> 5133.4998
$ws.Cells.Item(2, 12).Value = 3101988.36  # L2: 2929677.36 -> 3101988.36
$ws.Cells.Item(2, 13).Value = -5020.4998  # M2: -5482.8181 -> -5020.4998
$ws.Cells.Item(2, 14).Value = -3102214.36  # N2: -2929903.36 -> -3102214.36
$ws.Cells.Item(25, 8).Value = 88572.89  # H25: 2626 -> 88572.89
$ws.Cells.Item(25, 10).Value = 99625.75  # J25: 3245 -> 99625.75
$ws.Cells.Item(25, 12).Value = 298877.25  # L25: 9735 -> 298877.25
$ws.Cells.Item(25, 14).Value = -299215.25  # N25: -10073 -> -299215.25
$ws.Cells.Item(30, 8).Value = 88572.89  # H30: 2626 -> 88572.89
$ws.Cells.Item(30, 10).Value = 99625.75  # J30: 3245 -> 99625.75
$ws.Cells.Item(30, 12).Value = 298877.25  # L30: 9735 -> 298877.25
$ws.Cells.Item(30, 14).Value = -299081.25  # N30: -9939 -> -299081.25
$ws.Cells.Item(119, 8).Value = 4500  # H119: 3073.111 -> 4500
$ws.Cells.Item(119, 9).Value = 5750  # I119: 2457.25 -> 5750
$ws.Cells.Item(119, 10).Value = 2000  # J119: 8000 -> 2000
$ws.Cells.Item(119, 11).Value = 17250  # K119: 7371.75 -> 17250
$ws.Cells.Item(119, 12).Value = 6000  # L119: 24000 -> 6000
$ws.Cells.Item(119, 13).Value = -12412  # M119: -2533.75 -> -12412
$ws.Cells.Item(119, 14).Value = -15676  # N119: -33676 -> -15676
$ws.Cells.Item(132, 8).Value = 1472967.6  # H132: 835147.6 -> 1472967.6
$ws.Cells.Item(132, 9).Value = 2737.5  # I132: 2158.524 -> 2737.5
$ws.Cells.Item(132, 10).Value = 5001520  # J132: 2778789 -> 5001520
$ws.Cells.Item(132, 11).Value = 24637.5  # K132: 19426.716 -> 24637.5
$ws.Cells.Item(132, 12).Value = 45013680  # L132: 25009101 -> 45013680
$ws.Cells.Item(132, 13).Value = -22107.5  # M132: -16896.716 -> -22107.5
$ws.Cells.Item(132, 14).Value = -45018740  # N132: -25014161 -> -45018740
$ws.Cells.Item(134, 8).Value = 9603.75  # H134: 9987.691999999999 -> 9603.75
$ws.Cells.Item(134, 9).Value = 9207.5  # I134: 10643.333 -> 9207.5
$ws.Cells.Item(134, 10).Value = 10000  # J134: 9425.714 -> 10000
$ws.Cells.Item(134, 11).Value = 27622.5  # K134: 31929.999 -> 27622.5
$ws.Cells.Item(134, 12).Value = 30000  # L134: 28277.142 -> 30000
$ws.Cells.Item(134, 13).Value = -22552.5  # M134: -26859.999 -> -22552.5
$ws.Cells.Item(134, 14).Value = -40140  # N134: -38417.142 -> -40140
$ws.Cells.Item(139, 8).Value = 1467.5238  # H139: 1874.8948 -> 1467.5238
$ws.Cells.Item(139, 9).Value = 772.7143  # I139: 1235.4546 -> 772.7143
$ws.Cells.Item(139, 10).Value = 2857.1428  # J139: 2754.125 -> 2857.1428
$ws.Cells.Item(139, 11).Value = 2318.1429  # K139: 3706.3638 -> 2318.1429
$ws.Cells.Item(139, 12).Value = 8571.428400000001  # L139: 8262.375 -> 8571.428400000001
$ws.Cells.Item(139, 13).Value = 2821.8571  # M139: 1433.6362 -> 2821.8571
$ws.Cells.Item(139, 14).Value = -18851.4284  # N139: -18542.375 -> -18851.4284
$ws.Cells.Item(140, 8).Value = 2667.7693  # H140: 2816.5454 -> 2667.7693
$ws.Cells.Item(140, 9).Value = 852.8182  # I140: 898.2 -> 852.8182
$ws.Cells.Item(140, 10).Value = 12650  # J140: 22000 -> 12650
$ws.Cells.Item(140, 11).Value = 2558.4546  # K140: 2694.6 -> 2558.4546
$ws.Cells.Item(140, 12).Value = 37950  # L140: 66000 -> 37950
$ws.Cells.Item(140, 13).Value = 2621.5454  # M140: 2485.4 -> 2621.5454
$ws.Cells.Item(140, 14).Value = -48310  # N140: -76360 -> -48310

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(51, 8).Value = 27559.8  # H51: 27933.166 -> 27559.8
$ws.Cells.Item(51, 10).Value = 27559.8  # J51: 27933.166 -> 27559.8
$ws.Cells.Item(51, 12).Value = 27559.8  # L51: 27933.166 -> 27559.8
$ws.Cells.Item(51, 14).Value = -28577.8  # N51: -28951.166 -> -28577.8
$ws.Cells.Item(57, 8).Value = 15518  # H57: 15058.167 -> 15518
$ws.Cells.Item(57, 10).Value = 16919.8  # J57: 16290.728 -> 16919.8
$ws.Cells.Item(57, 12).Value = 16919.8  # L57: 16290.728 -> 16919.8
$ws.Cells.Item(57, 14).Value = -18559.8  # N57: -17930.728 -> -18559.8
$ws.Cells.Item(132, 8).Value = 2634186.5  # H132: 2859885 -> 2634186.5
$ws.Cells.Item(132, 9).Value = 2767.8  # I132: 2927.6086 -> 2767.8
$ws.Cells.Item(132, 10).Value = 7694607  # J132: 8335720 -> 7694607
$ws.Cells.Item(132, 11).Value = 8303.400000000001  # K132: 8782.825800000001 -> 8303.400000000001
$ws.Cells.Item(132, 12).Value = 23083821  # L132: 25007160 -> 23083821
$ws.Cells.Item(132, 13).Value = -5773.400000000001  # M132: -6252.825800000001 -> -5773.400000000001
$ws.Cells.Item(132, 14).Value = -23088881  # N132: -25012220 -> -23088881

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(40, 8).Value = 2376  # H40: 2700.5 -> 2376
$ws.Cells.Item(40, 9).Value = 2286.8572  # I40: 2657.7144 -> 2286.8572
$ws.Cells.Item(40, 11).Value = 2286.8572  # K40: 2657.7144 -> 2286.8572
$ws.Cells.Item(40, 13).Value = -2150.8572  # M40: -2521.7144 -> -2150.8572
$ws.Cells.Item(42, 8).Value = 20000  # H42: 0 -> 20000
$ws.Cells.Item(42, 10).Value = 20000  # J42: 0 -> 20000
$ws.Cells.Item(42, 12).Value = 20000  # L42: 0 -> 20000
$ws.Cells.Item(42, 14).Value = -21126  # N42: None -> -21126
$ws.Cells.Item(49, 8).Value = 20000  # H49: 0 -> 20000
$ws.Cells.Item(49, 10).Value = 20000  # J49: 0 -> 20000
$ws.Cells.Item(49, 12).Value = 20000  # L49: 0 -> 20000
$ws.Cells.Item(49, 14).Value = -20294  # N49: None -> -20294
$ws.Cells.Item(61, 8).Value = 1939.5312  # H61: 2374.6365 -> 1939.5312
$ws.Cells.Item(61, 9).Value = 1717.9615  # I61: 2256.2666 -> 1717.9615
$ws.Cells.Item(61, 10).Value = 2899.6667  # J61: 2628.2856 -> 2899.6667
$ws.Cells.Item(61, 11).Value = 1717.9615  # K61: 2256.2666 -> 1717.9615
$ws.Cells.Item(61, 12).Value = 2899.6667  # L61: 2628.2856 -> 2899.6667
$ws.Cells.Item(61, 13).Value = -1515.9615  # M61: -2054.2666 -> -1515.9615
$ws.Cells.Item(61, 14).Value = -3303.6667  # N61: -3032.2856 -> -3303.6667
$ws.Cells.Item(93, 8).Value = 51705.93  # H93: 51804.715 -> 51705.93
$ws.Cells.Item(93, 9).Value = 1198.9  # I93: 1337.2 -> 1198.9
$ws.Cells.Item(93, 11).Value = 1198.9  # K93: 1337.2 -> 1198.9
$ws.Cells.Item(93, 13).Value = 49.09999999999991  # M93: -89.20000000000005 -> 49.09999999999991
$ws.Cells.Item(100, 8).Value = 1391.6666  # H100: 1837.5 -> 1391.6666
$ws.Cells.Item(100, 9).Value = 1280  # I100: 1725 -> 1280
$ws.Cells.Item(100, 11).Value = 1280  # K100: 1725 -> 1280
$ws.Cells.Item(100, 13).Value = -739  # M100: -1184 -> -739
$ws.Cells.Item(113, 8).Value = 1939.5312  # H113: 2374.6365 -> 1939.5312
$ws.Cells.Item(113, 9).Value = 1717.9615  # I113: 2256.2666 -> 1717.9615
$ws.Cells.Item(113, 10).Value = 2899.6667  # J113: 2628.2856 -> 2899.6667
$ws.Cells.Item(113, 11).Value = 1717.9615  # K113: 2256.2666 -> 1717.9615
$ws.Cells.Item(113, 12).Value = 2899.6667  # L113: 2628.2856 -> 2899.6667
$ws.Cells.Item(113, 13).Value = 452.0385000000001  # M113: -86.26659999999993 -> 452.0385000000001
$ws.Cells.Item(113, 14).Value = -7239.6667  # N113: -6968.2856 -> -7239.6667

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(113, 8).Value = 577.46155  # H113: 723.7895 -> 577.46155
$ws.Cells.Item(113, 9).Value = 445.7  # I113: 588.1539 -> 445.7
$ws.Cells.Item(113, 10).Value = 1016.6667  # J113: 1017.6667 -> 1016.6667
$ws.Cells.Item(113, 11).Value = 1337.1  # K113: 1764.4617 -> 1337.1
$ws.Cells.Item(113, 12).Value = 3050.0001  # L113: 3053.0001 -> 3050.0001
$ws.Cells.Item(113, 13).Value = 832.9000000000001  # M113: 405.5382999999999 -> 832.9000000000001
$ws.Cells.Item(113, 14).Value = -7390.0001  # N113: -7393.0001 -> -7390.0001
$ws.Cells.Item(136, 8).Value = 3699.9355  # H136: 2027.5079 -> 3699.9355
$ws.Cells.Item(136, 9).Value = 4130.5625  # I136: 1679.234 -> 4130.5625
$ws.Cells.Item(136, 10).Value = 3240.6  # J136: 3050.5625 -> 3240.6
$ws.Cells.Item(136, 11).Value = 12391.6875  # K136: 5037.701999999999 -> 12391.6875
$ws.Cells.Item(136, 12).Value = 9721.799999999999  # L136: 9151.6875 -> 9721.799999999999
$ws.Cells.Item(136, 13).Value = -9841.6875  # M136: -2487.701999999999 -> -9841.6875
$ws.Cells.Item(136, 14).Value = -14821.8  # N136: -14251.6875 -> -14821.8
